$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell AC1 with "No Of Months" (no special style, matches default)
$ws.Range("AC1").Value = "No Of Months"

# Set width for new column AC (29). COM ColumnWidth is offset from the raw
# OOXML column width by the standard 5px gridline padding (~0.8333 chars at
# the workbook's default font), so subtract that to land on width="13".
$ws.Columns.Item(29).ColumnWidth = 12.166666666666666

# Update view: scroll so AB1 is the top-left cell, and select AJ12
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 28
$ws.Range("AJ12").Select()
